$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R: header "backup" ---
# Copy the format from Q1 (bold, bordered, centered header style) to R1
$ws.Cells.Item(1, 17).Copy()
$ws.Cells.Item(1, 18).PasteSpecial(-4122)
$ws.Cells.Item(1, 18).Value2 = "backup"

# --- Fill column R with 0 for existing data rows 2-75 ---
for ($r = 2; $r -le 75; $r++) {
    $ws.Cells.Item($r, 18).Value2 = 0
}

# --- Correct Q column values that were recalculated to 0 ---
$qZeroRows = @(28, 33, 36, 40, 47, 49)
foreach ($r in $qZeroRows) {
    $ws.Cells.Item($r, 17).Value2 = 0
}

# --- Append new rows 76-81 ---
# Copy number-format (date style) from A75 for the new date cells in column A
$newRows = @(
    @(45474, 5255.534415453605, 5659.173420131944, 4497.620533908113, 4909.33740234375, 69638843, 2024, 7, 1, 0, 0, 0, 27, 1, 0, 0),
    @(45505, 4908.489395941672, 4933.320949500038, 4467.802483157853, 4667.10400390625, 39407162, 2024, 8, 1, 0, 0, 0, 31, 0, 0, 0),
    @(45536, 4699, 4950, 4172, 4420.64990234375, 38452011, 2024, 9, 1, 0, 0, 0, 35, 0, 0, 0),
    @(45566, 4428, 4676.60009765625, 4075.199951171875, 4246.7001953125, 32512947, 2024, 10, 1, 0, 0, 0, 40, 0, 0, 0),
    @(45597, 4269.7998046875, 4529.5, 3920.35009765625, 4476.85009765625, 29976034, 2024, 11, 1, 0, 0, 0, 44, 0, 0, 2),
    @(45627, 4475, 4755, 4050, 4081.050048828125, 27164344, 2024, 12, 1, 0, 0, 0, 48, 0, 0, 0)
)

$rowNum = 76
foreach ($row in $newRows) {
    # Column A: Datetime (copy date number format/style from A75)
    $ws.Cells.Item(75, 1).Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowNum, 1).Value2 = $row[0]

    # Column B: Open
    $ws.Cells.Item($rowNum, 2).Value2 = $row[1]
    # Column C: High
    $ws.Cells.Item($rowNum, 3).Value2 = $row[2]
    # Column D: Low
    $ws.Cells.Item($rowNum, 4).Value2 = $row[3]
    # Column E: Close
    $ws.Cells.Item($rowNum, 5).Value2 = $row[4]
    # Column F: Adj Close -- left blank (no data yet)
    $ws.Cells.Item($rowNum, 6).Value2 = ""
    # Column G: Volume
    $ws.Cells.Item($rowNum, 7).Value2 = $row[5]
    # Column H: Year
    $ws.Cells.Item($rowNum, 8).Value2 = $row[6]
    # Column I: Month
    $ws.Cells.Item($rowNum, 9).Value2 = $row[7]
    # Column J: Day
    $ws.Cells.Item($rowNum, 10).Value2 = $row[8]
    # Column K: Hour
    $ws.Cells.Item($rowNum, 11).Value2 = $row[9]
    # Column L: Minute
    $ws.Cells.Item($rowNum, 12).Value2 = $row[10]
    # Column M: Second
    $ws.Cells.Item($rowNum, 13).Value2 = $row[11]
    # Column N: Week
    $ws.Cells.Item($rowNum, 14).Value2 = $row[12]
    # Column O: isPivot
    $ws.Cells.Item($rowNum, 15).Value2 = $row[13]
    # Column P: two_line_structure
    $ws.Cells.Item($rowNum, 16).Value2 = $row[14]
    # Column Q: detect_structure
    $ws.Cells.Item($rowNum, 17).Value2 = $row[15]
    # Column R: backup -- left blank (no data yet)
    $ws.Cells.Item($rowNum, 18).Value2 = ""

    $rowNum++
}
